$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the first name in row 3 (Branka -> Branko)
$ws.Range("A3").Value = "Branko"

# Match the saved selection state from the edit (active cell moved to A3)
$ws.Range("A3").Select()
